$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text value looks like a plain number (e.g. "210.47") need to be
# forced to stay text - otherwise Excel auto-converts them to a numeric cell.
# We flip the cell to Text format only for the duration of the write, then restore
# the default "Normal" style so no stray formatting is left behind.

$ws.Range("D2").Value = '28.598.84'
$ws.Range("E2").Value = '  +0.64%  '
$ws.Range("D3").Value = '1.562.65'
$ws.Range("E3").Value = '  -0.76%  '
$ws.Range("E4").Value = '  -0.32%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '210.47'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.71%  '
$ws.Range("E6").Value = '  -0.99%  '
$ws.Range("E7").Value = '  -0.28%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '25.09'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +5.39%  '
$ws.Range("E9").Value = '  -0.87%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0586'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -0.25%  '
$ws.Range("E11").Value = '  +0.14%  '
$ws.Range("D12").Value = '1.786.79'
$ws.Range("E12").Value = '  -0.65%  '
$ws.Range("D13").Value = '1.562.96'
$ws.Range("E13").Value = '  -0.70%  '
$ws.Range("D14").Value = '28.607.05'
$ws.Range("E14").Value = '  +0.67%  '
$ws.Range("E15").Value = '  -0.21%  '
$ws.Range("E16").Value = '  -1.10%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '61.23'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -0.79%  '
$ws.Range("E18").Value = '  +0.38%  '
$ws.Range("E19").Value = '  -0.67%  '
$ws.Range("D20").Value = '0.0₃0678'
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.999'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -0.26%  '
$ws.Range("E22").Value = '  -1.20%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '9.00'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +0.33%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.09'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +1.45%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '151.01'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -0.02%  '
$ws.Range("E26").Value = '  -1.18%  '
$ws.Range("E27").Value = '  -0.08%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.00'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -0.21%  '
$ws.Range("E29").Value = '  -2.15%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.0462'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -4.01%  '
$ws.Range("E31").Value = '  -2.80%  '
$ws.Range("E32").Value = '  -0.44%  '
$ws.Range("D33").Value = '1.387.88'
$ws.Range("E33").Value = '  +0.38%  '
$ws.Range("E34").Value = '  -4.33%  '
$ws.Range("E35").Value = '  -3.77%  '
$ws.Range("E36").Value = '  -1.49%  '
$ws.Range("E37").Value = '  +1.73%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.30'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -2.22%  '
$ws.Range("E39").Value = '  -1.08%  '
$ws.Range("E40").Value = '  +1.94%  '
$ws.Range("E41").Value = '  -0.47%  '
$ws.Range("E42").Value = '  -0.24%  '
$ws.Range("E43").Value = '  -1.88%  '
$ws.Range("E44").Value = '  -2.28%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '63.97'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +2.47%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '5.23'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -2.31%  '
$ws.Range("D47").Value = '1.699.13'
$ws.Range("E47").Value = '  -0.64%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.870'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -5.27%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '85.19'
$ws.Range("D49").Style = "Normal"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '43.22'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +6.62%  '
$ws.Range("D51").Value = '0.0₆0101'
$ws.Range("E51").Value = '  +0.02%  '
